# Change all the database and xlsx data schema to lower case (Postgres-friendly).
$wb = $excel.ActiveWorkbook

# --- Rename sheets to snake_case / lower case ---
$wb.Worksheets.Item("Trade").Name = "trade"
$wb.Worksheets.Item("WorkMethod").Name = "work_method"
$wb.Worksheets.Item("WorkMethodDependency").Name = "work_method_dependency"
$wb.Worksheets.Item("Space").Name = "space"
$wb.Worksheets.Item("Design").Name = "design"

# --- trade: header Name -> name ---
$ws = $wb.Worksheets.Item("trade")
$ws.Range("A1").Value = "name"

# --- work_method: headers TradeName/Method/Productivity -> trade_name/method/productivity ---
$ws = $wb.Worksheets.Item("work_method")
$ws.Range("A1").Value = "trade_name"
$ws.Range("B1").Value = "method"
$ws.Range("C1").Value = "productivity"

# --- work_method_dependency: headers Predecessor/Successor -> predecessor/successor ---
$ws = $wb.Worksheets.Item("work_method_dependency")
$ws.Range("A1").Value = "predecessor"
$ws.Range("B1").Value = "successor"

# --- space: collapse WorkSpace(id)/Name(floor) table into a single name column ---
$ws = $wb.Worksheets.Item("space")
$ws.Range("A1:B6").ClearContents()
$ws.Range("A1").Value = "name"
$ws.Range("A2").Value = "floor_1"
$ws.Range("A3").Value = "floor_2"

# --- design: headers WorkMethod/WorkSpace/Quantity -> work_method/space/quantity ---
$ws = $wb.Worksheets.Item("design")
$ws.Range("A1").Value = "work_method"
$ws.Range("B1").Value = "space"
$ws.Range("C1").Value = "quantity"

# --- match final view state: "space" tab active with column A selected ---
$ws = $wb.Worksheets.Item("space")
$ws.Activate() | Out-Null
$ws.Columns("A").Select() | Out-Null
